$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("X1").Value = "user_mobile"

for ($row = 2; $row -le 23; $row++) {
    $ws.Cells.Item($row, 24).Value = 7896302536
}
